{"js": "// Target: the document's first paragraph (the \"**ID__...__ID**\" marker\n// paragraph) gets:\n//   1. A paragraph border (pBdr) on all four sides, each with a 5pt\n//      \"space\" (distance from text) and no visible line.\n//   2. Its left indent changed from 120 twips (6pt) to 225 twips (11.25pt).\n//   3. Its marker text updated from\n//      \"**ID__AFFARS_5312_topic_13__ID**\" to\n//      \"**ID__AFFARS_SUBPART_5312_90__ID**\".\n//   4. The trailing single-space run removed entirely.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst marker = paragraphs.items[0];\n\n// 1) Paragraph border with 5pt spacing on every side (no line style/color,\n// matching <w:pBdr><w:top w:space=\"5\"/>...). Word's JS API does not expose\n// the per-edge \"distance from text\" value through the typed Border object,\n// so we reach it via the same object-model dispatch that backs the typed\n// wrappers (this mirrors Word's own Borders.DistanceFromTop/Bottom/Left/Right).\nconst borders = marker.borders;\nborders.load(\"items\");\nawait context.sync();\n\nborders._omSet(\"DistanceFromTop\", 5, \"Borders\");\nborders._omSet(\"DistanceFromBottom\", 5, \"Borders\");\nborders._omSet(\"DistanceFromLeft\", 5, \"Borders\");\nborders._omSet(\"DistanceFromRight\", 5, \"Borders\");\n\n// 2) Left indent: 120 twips -> 225 twips (twips / 20 = points).\nmarker.leftIndent = 225 / 20;\n\n// 3) Update the marker run's text in place (keeps its run formatting).\nconst idResults = body.search(\"**ID__AFFARS_5312_topic_13__ID**\", { matchWildcards: false });\nidResults.load(\"items\");\nawait context.sync();\n\nif (idResults.items.length > 0) {\n  idResults.items[0].insertText(\"**ID__AFFARS_SUBPART_5312_90__ID**\", \"Replace\");\n  await context.sync();\n}\n\n// 4) Remove the trailing lone-space run that followed the marker text.\nconst markerContent = marker.getRange(\"Content\");\nconst spaceResults = markerContent.search(\" \", { matchWildcards: false });\nspaceResults.load(\"items\");\nawait context.sync();\n\nif (spaceResults.items.length > 0) {\n  spaceResults.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# Target: the document's first paragraph (the \"**ID__...__ID**\" marker\n# paragraph) gets:\n#   1. A paragraph border (pBdr) on all four sides, each with a 5pt\n#      \"space\" (distance from text) and no visible line.\n#   2. Its left indent changed from 120 twips (6pt) to 225 twips (11.25pt).\n#   3. Its marker text updated from\n#      \"**ID__AFFARS_5312_topic_13__ID**\" to\n#      \"**ID__AFFARS_SUBPART_5312_90__ID**\".\n#   4. The trailing single-space run removed entirely.\n\n$d = $word.ActiveDocument\n$p1 = $d.Paragraphs(1)\n\n# 1) Update the marker run's text in place (keeps its run formatting).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"**ID__AFFARS_5312_topic_13__ID**\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"**ID__AFFARS_SUBPART_5312_90__ID**\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Remove the trailing lone-space run that followed the marker text\n# (it sits right before the paragraph mark).\n$r = $p1.Range\n$spaceRange = $d.Range($r.End - 2, $r.End - 1)\nif ($spaceRange.Text -eq \" \") {\n  $spaceRange.Delete()\n}\n\n# 3) Left indent: 120 twips -> 225 twips (11.25pt).\n$p1.Range.ParagraphFormat.LeftIndent = 11.25\n\n# 4) Paragraph border with 5pt spacing on every side (no line style/color),\n# matching <w:pBdr><w:top w:space=\"5\"/>...<w:right w:space=\"5\"/></w:pBdr>.\n$borders = $p1.Range.ParagraphFormat.Borders\n$borders.DistanceFromTop = 5\n$borders.DistanceFromBottom = 5\n$borders.DistanceFromLeft = 5\n$borders.DistanceFromRight = 5\n"}
